# Rename the "Evidence" column header (N1) to "Preuve"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "Preuve"

# Move the active selection from G19 to N2, matching the saved cursor position
$ws.Range("N2").Select()
